$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7732.6665
$ws.Range("I51").Value = 6375.5
$ws.Range("J51").Value = 8120.4287
$ws.Range("K51").Value = 6375.5
$ws.Range("L51").Value = 8120.4287
$ws.Range("M51").Value = -5891.5
$ws.Range("N51").Value = -9088.4287
$ws.Range("H100").Value = 3801
$ws.Range("I100").Value = 4600
$ws.Range("K100").Value = 4600
$ws.Range("M100").Value = -4059
$ws.Range("H106").Value = 1487.5
$ws.Range("J106").Value = 875
$ws.Range("L106").Value = 875
$ws.Range("N106").Value = -2137
$ws.Range("H137").Value = 12508.484
$ws.Range("I137").Value = 3574.3572
$ws.Range("K137").Value = 10723.0716
$ws.Range("M137").Value = -8173.071599999999
$ws.Range("H138").Value = 2556.5483
$ws.Range("J138").Value = 3229.75
$ws.Range("L138").Value = 9689.25
$ws.Range("N138").Value = -19969.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 106.22222
$ws.Range("I4").Value = 132
$ws.Range("K4").Value = 132
$ws.Range("M4").Value = -16
$ws.Range("H61").Value = 10948.546
$ws.Range("I61").Value = 2959.7
$ws.Range("J61").Value = 14421.956
$ws.Range("K61").Value = 2959.7
$ws.Range("L61").Value = 14421.956
$ws.Range("M61").Value = -2747.7
$ws.Range("N61").Value = -14845.956
$ws.Range("H136").Value = 10948.546
$ws.Range("I136").Value = 2959.7
$ws.Range("J136").Value = 14421.956
$ws.Range("K136").Value = 8879.099999999999
$ws.Range("L136").Value = 43265.868
$ws.Range("M136").Value = -6329.099999999999
$ws.Range("N136").Value = -48365.868
$ws.Range("H137").Value = 62261.816
$ws.Range("J137").Value = 63110
$ws.Range("L137").Value = 63110
$ws.Range("N137").Value = -73310
$ws.Range("H139").Value = 63342.57
$ws.Range("J139").Value = 64753.617
$ws.Range("L139").Value = 64753.617
$ws.Range("N139").Value = -75033.617

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 2354.6667
$ws.Range("I25").Value = 376
$ws.Range("J25").Value = 4333.3335
$ws.Range("K25").Value = 376
$ws.Range("L25").Value = 4333.3335
$ws.Range("M25").Value = -141
$ws.Range("N25").Value = -4803.3335
$ws.Range("H51").Value = 57800
$ws.Range("J51").Value = 57800
$ws.Range("L51").Value = 57800
$ws.Range("N51").Value = -58782
$ws.Range("H92").Value = 58993
$ws.Range("J92").Value = 58993
$ws.Range("L92").Value = 58993
$ws.Range("N92").Value = -63985
$ws.Range("H94").Value = 11002
$ws.Range("I94").Value = 4504
$ws.Range("K94").Value = 4504
$ws.Range("M94").Value = -4053
$ws.Range("H134").Value = 15930.087
$ws.Range("J134").Value = 26657.666
$ws.Range("L134").Value = 79972.99800000001
$ws.Range("N134").Value = -85042.99800000001
$ws.Range("H141").Value = 138999
$ws.Range("J141").Value = 138999
$ws.Range("L141").Value = 138999
$ws.Range("N141").Value = -149359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 861.05554
$ws.Range("I22").Value = 381.44446
$ws.Range("J22").Value = 1340.6666
$ws.Range("K22").Value = 381.44446
$ws.Range("L22").Value = 1340.6666
$ws.Range("M22").Value = -31.44445999999999
$ws.Range("N22").Value = -2040.6666
$ws.Range("H134").Value = 37044972
$ws.Range("I134").Value = 2401.5454
$ws.Range("J134").Value = 62511736
$ws.Range("K134").Value = 7204.6362
$ws.Range("L134").Value = 187535208
$ws.Range("M134").Value = -4669.6362
$ws.Range("N134").Value = -187540278
$ws.Range("H141").Value = 281790.75
$ws.Range("J141").Value = 281790.75
$ws.Range("L141").Value = 281790.75
$ws.Range("N141").Value = -292150.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 124.1
$ws.Range("I23").Value = 125.5
$ws.Range("J23").Value = 123.75
$ws.Range("K23").Value = 376.5
$ws.Range("L23").Value = 371.25
$ws.Range("M23").Value = -141.5
$ws.Range("N23").Value = -841.25
$ws.Range("H117").Value = 2523.25
$ws.Range("I117").Value = 2750
$ws.Range("J117").Value = 2477.9
$ws.Range("K117").Value = 8250
$ws.Range("L117").Value = 7433.700000000001
$ws.Range("M117").Value = -4808
$ws.Range("N117").Value = -14317.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 17037.25
$ws.Range("I80").Value = 20390.8
$ws.Range("J80").Value = 14641.857
$ws.Range("K80").Value = 20390.8
$ws.Range("L80").Value = 14641.857
$ws.Range("M80").Value = -19392.8
$ws.Range("N80").Value = -16637.857
$ws.Range("H83").Value = 17037.25
$ws.Range("I83").Value = 20390.8
$ws.Range("J83").Value = 14641.857
$ws.Range("K83").Value = 101954
$ws.Range("L83").Value = 73209.285
$ws.Range("M83").Value = -96962
$ws.Range("N83").Value = -83193.285
$ws.Range("H97").Value = 1366.238
$ws.Range("I97").Value = 1553
$ws.Range("K97").Value = 1553
$ws.Range("M97").Value = -1057
$ws.Range("H126").Value = 4046.6584
$ws.Range("I126").Value = 3227.6785
$ws.Range("K126").Value = 9683.0355
$ws.Range("M126").Value = -7213.0355
$ws.Range("H132").Value = 6666.483
$ws.Range("I132").Value = 2260.55
$ws.Range("J132").Value = 16457.445
$ws.Range("K132").Value = 6781.650000000001
$ws.Range("L132").Value = 49372.335
$ws.Range("M132").Value = -4251.650000000001
$ws.Range("N132").Value = -54432.335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -137
$ws.Range("H21").Value = 1902.5
$ws.Range("I21").Value = 2127.25
$ws.Range("K21").Value = 2127.25
$ws.Range("M21").Value = -1953.25
$ws.Range("H28").Value = 250
$ws.Range("I28").Value = 250
$ws.Range("K28").Value = 250
$ws.Range("M28").Value = -18
$ws.Range("H35").Value = 1255.3334
$ws.Range("I35").Value = 1255.3334
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1255.3334
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -919.3334
$ws.Range("N35").ClearContents() | Out-Null
$ws.Range("H37").Value = 250
$ws.Range("I37").Value = 250
$ws.Range("K37").Value = 250
$ws.Range("M37").Value = -143
$ws.Range("H40").Value = 4814.5356
$ws.Range("J40").Value = 5462.727
$ws.Range("L40").Value = 5462.727
$ws.Range("N40").Value = -5734.727
$ws.Range("H122").Value = 5791.7144
$ws.Range("I122").Value = 5004.316
$ws.Range("K122").Value = 15012.948
$ws.Range("M122").Value = -12562.948
$ws.Range("H127").Value = 286999.5
$ws.Range("I127").Value = 70000
$ws.Range("J127").Value = 359332.66
$ws.Range("K127").Value = 70000
$ws.Range("L127").Value = 359332.66
$ws.Range("M127").Value = -65040
$ws.Range("N127").Value = -369252.66
$ws.Range("H132").Value = 1153340.6
$ws.Range("J132").Value = 2015340.5
$ws.Range("L132").Value = 6046021.5
$ws.Range("N132").Value = -6051081.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 98641.14
$ws.Range("J46").Value = 98641.14
$ws.Range("L46").Value = 98641.14
$ws.Range("N46").Value = -99103.14
$ws.Range("H48").Value = 16000
$ws.Range("I48").Value = 13400
$ws.Range("J48").Value = 29000
$ws.Range("K48").Value = 13400
$ws.Range("L48").Value = 29000
$ws.Range("M48").Value = -12831
$ws.Range("N48").Value = -30138
$ws.Range("H122").Value = 3678.3076
$ws.Range("I122").Value = 2378.8
$ws.Range("K122").Value = 7136.400000000001
$ws.Range("M122").Value = -4686.400000000001
$ws.Range("H132").Value = 6540
$ws.Range("I132").Value = 1986.762
$ws.Range("J132").Value = 15232.546
$ws.Range("K132").Value = 5960.286
$ws.Range("L132").Value = 45697.638
$ws.Range("M132").Value = -3430.286
$ws.Range("N132").Value = -50757.638
$ws.Range("H134").Value = 98641.14
$ws.Range("J134").Value = 98641.14
$ws.Range("L134").Value = 295923.42
$ws.Range("N134").Value = -300993.42

